$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Bump up a few row heights (28.8pt -> 30pt), and give row 3 the same
# explicit 30pt height.
$ws.Rows.Item(3).RowHeight = 30
$ws.Rows.Item(4).RowHeight = 30
$ws.Rows.Item(6).RowHeight = 30
$ws.Rows.Item(11).RowHeight = 30

# Fill in the "finished date" column for the last three rows, copying the
# date formatting already used in column B so the new cells match.
$ws.Range("C9").Value2 = 44761
$ws.Range("B9").Copy()
$ws.Range("C9").PasteSpecial(-4122)

$ws.Range("C10").Value2 = 44761
$ws.Range("B10").Copy()
$ws.Range("C10").PasteSpecial(-4122)

$ws.Range("C11").Value2 = 44761
$ws.Range("B11").Copy()
$ws.Range("C11").PasteSpecial(-4122)

$excel.CutCopyMode = 0

# Move the selection to the last filled-in cell, matching where the user
# ended up after entering the new data.
[void]$ws.Range("C11").Select()
